$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Uncheck every "[ X]" checkbox in the checklist -> "[  ]" (the X becomes
#    a plain space, same as someone selecting the X and typing a space).
#    A single Find/Replace (Replace:=wdReplaceAll) over the whole document
#    mirrors exactly what Word's Find & Replace does, regardless of how the
#    literal "[ X]" text happens to be split across runs.
# ---------------------------------------------------------------------------
$rngAll = $d.Content
$rngAll.Find.Execute("[ X]", $true, $false, $false, $false, $false, $true, 1, $false, "[  ]", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark. It used to sit right after "DATE: " ; Word
#    re-drops it at the location of the most recent edit, which here is the
#    very first checkbox (the one in front of "...revision levels and
#    titles..."), landing immediately before the closing "]".
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$firstBox = $d.Range(0, 400)
$firstBox.Find.Execute("[  ]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $firstBox.End - 1

$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
